$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 15627388
$ws.Range("I113").Value = 2525
$ws.Range("J113").Value = 31252250
$ws.Range("K113").Value = 2525
$ws.Range("L113").Value = 31252250
$ws.Range("M113").Value = 729
$ws.Range("N113").Value = -31258758
$ws.Range("H128").Value = 79800
$ws.Range("J128").Value = 79800
$ws.Range("L128").Value = 79800
$ws.Range("N128").Value = -89760
$ws.Range("H137").Value = 1521
$ws.Range("I137").Value = 1111.2128
$ws.Range("J137").Value = 3271.9092
$ws.Range("K137").Value = 3333.6384
$ws.Range("L137").Value = 9815.7276
$ws.Range("M137").Value = -783.6383999999998
$ws.Range("N137").Value = -14915.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 659628.5
$ws.Range("I122").Value = 856760.1
$ws.Range("J122").Value = 2523.111
$ws.Range("K122").Value = 2570280.3
$ws.Range("L122").Value = 7569.333
$ws.Range("M122").Value = -2567830.3
$ws.Range("N122").Value = -12469.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3011
$ws.Range("I107").Value = 3011
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3011
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1091
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 40230.07
$ws.Range("I134").Value = 6770.8184
$ws.Range("J134").Value = 145387.72
$ws.Range("K134").Value = 20312.4552
$ws.Range("L134").Value = 436163.16
$ws.Range("M134").Value = -17777.4552
$ws.Range("N134").Value = -441233.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 19001.875
$ws.Range("J21").Value = 19001.875
$ws.Range("L21").Value = 19001.875
$ws.Range("N21").Value = -19471.875
$ws.Range("H23").Value = 250009000
$ws.Range("J23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12480
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 250009000
$ws.Range("J27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("N27").Value = -12384
$ws.Range("H41").Value = 19500
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19500
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19500
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -20356
$ws.Range("H42").Value = 14000
$ws.Range("I42").Value = 14000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 14000
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("M42").Value = -13407
$ws.Range("H56").Value = 28000
$ws.Range("I56").Value = 28000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 28000
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("M56").Value = -27155
$ws.Range("H99").Value = 55000
$ws.Range("I99").Value = 55000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 55000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -53502
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1034.2188
$ws.Range("I107").Value = 1101.6842
$ws.Range("J107").Value = 935.61536
$ws.Range("K107").Value = 1101.6842
$ws.Range("L107").Value = 935.61536
$ws.Range("M107").Value = 818.3158000000001
$ws.Range("N107").Value = -4775.61536
$ws.Range("H118").Value = 39999.94
$ws.Range("J118").Value = 39999.94
$ws.Range("L118").Value = 39999.94
$ws.Range("N118").Value = -43313.94
$ws.Range("H126").Value = 55000
$ws.Range("I126").Value = 55000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 165000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -162530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5078.5
$ws.Range("I70").Value = 4849.9
$ws.Range("J70").Value = 5650
$ws.Range("K70").Value = 4849.9
$ws.Range("L70").Value = 5650
$ws.Range("M70").Value = -4579.9
$ws.Range("N70").Value = -6190
$ws.Range("H73").Value = 5078.5
$ws.Range("I73").Value = 4849.9
$ws.Range("J73").Value = 5650
$ws.Range("K73").Value = 4849.9
$ws.Range("L73").Value = 5650
$ws.Range("M73").Value = -3913.9
$ws.Range("N73").Value = -7522
$ws.Range("H122").Value = 39441084
$ws.Range("I122").Value = 59158796
$ws.Range("J122").Value = 5656
$ws.Range("K122").Value = 177476388
$ws.Range("L122").Value = 16968
$ws.Range("M122").Value = -177473938
$ws.Range("N122").Value = -21868
$ws.Range("H123").Value = 18599.727
$ws.Range("J123").Value = 18599.727
$ws.Range("L123").Value = 18599.727
$ws.Range("N123").Value = -23499.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1640.4445
$ws.Range("I16").Value = 1770.25
$ws.Range("K16").Value = 1770.25
$ws.Range("M16").Value = -1600.25
$ws.Range("H61").Value = 2263.625
$ws.Range("I61").Value = 2040.8
$ws.Range("J61").Value = 2635
$ws.Range("K61").Value = 2040.8
$ws.Range("L61").Value = 2635
$ws.Range("M61").Value = -1838.8
$ws.Range("N61").Value = -3039
$ws.Range("H113").Value = 2263.625
$ws.Range("I113").Value = 2040.8
$ws.Range("J113").Value = 2635
$ws.Range("K113").Value = 2040.8
$ws.Range("L113").Value = 2635
$ws.Range("M113").Value = 129.2
$ws.Range("N113").Value = -6975
$ws.Range("H122").Value = 8939271
$ws.Range("I122").Value = 8939271
$ws.Range("K122").Value = 26817813
$ws.Range("M122").Value = -26815363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7685.7144
$ws.Range("J54").Value = 7685.7144
$ws.Range("L54").Value = 7685.7144
$ws.Range("N54").Value = -8725.7144
$ws.Range("H86").Value = 20012
$ws.Range("J86").Value = 20012
$ws.Range("L86").Value = 20012
$ws.Range("N86").Value = -22258
$ws.Range("H89").Value = 20012
$ws.Range("J89").Value = 20012
$ws.Range("L89").Value = 100060
$ws.Range("N89").Value = -111292
$ws.Range("H122").Value = 1707.4375
$ws.Range("I122").Value = 1791.909
$ws.Range("J122").Value = 1521.6
$ws.Range("K122").Value = 5375.727000000001
$ws.Range("L122").Value = 4564.799999999999
$ws.Range("M122").Value = -2925.727000000001
$ws.Range("N122").Value = -9464.8
